# Apply the Ver-3.1.1 edit: update row 2 (Astronomy) values for the
# Betweenness value columns and rewrite the Closeness / Degree blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Betweenness block (columns B-G) : values only change ---
$ws.Range("C2").Value = 0.1983368870485608
$ws.Range("E2").Value = 0.1478279317385609
$ws.Range("G2").Value = 0.1410976472855235

# --- Closeness block (columns H-M) : nodes + values change ---
$ws.Range("H2").Value = "행성"
$ws.Range("I2").Value = 0.3500268749776042
$ws.Range("J2").Value = "태양계"
$ws.Range("K2").Value = 0.3492884216548666
$ws.Range("L2").Value = "케플러법칙"
$ws.Range("M2").Value = 0.3083104503992677

# --- Degree block (columns N-S) : node1/value1, node3/value3 and value2 change ---
$ws.Range("N2").Value = "태양계"
$ws.Range("O2").Value = 0.2779661016949153
$ws.Range("P2").Value = "행성"
$ws.Range("Q2").Value = 0.2711864406779661
$ws.Range("R2").Value = "케플러법칙"
$ws.Range("S2").Value = 0.2372881355932203
